# Dobrodosli izadajice na Mobu
# Adds two new test rows (C70803, C70804) to the QA test matrix on Sheet1,
# extends the AutoFilter / _FilterDatabase range to cover them, and moves
# the active selection the way the authoring tool left it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the AutoFilter range from A1:G24 to A1:G25 first (new last data
#     row is 25) while the sheet's used range still ends at row 25, so the
#     filter doesn't auto-expand to swallow the rows we are about to add. ---
$ws.AutoFilterMode = $false
$ws.Range("A1:G25").AutoFilter() | Out-Null

# --- Update the hidden _xlnm._FilterDatabase defined name to match ---
$names = $ws.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$25"
    }
}

# --- Append row 26: Current_Domestic_Accounts-Transactions-Filter_By_Amount_Invalid_[MOB_ANDROID] / C70803 ---
$ws.Range("A25:G25").Copy($ws.Range("A26:G26")) | Out-Null
$ws.Range("B26").Value = "Current_Domestic_Accounts-Transactions-Filter_By_Amount_Invalid_[MOB_ANDROID]"
$ws.Range("C26").Value = "C70803"

# --- Append row 27: Current_Domestic_Accounts-Transactions-List_[MOB_ANDROID] / C70804 ---
$ws.Range("A25:G25").Copy($ws.Range("A27:G27")) | Out-Null
$ws.Range("B27").Value = "Current_Domestic_Accounts-Transactions-List_[MOB_ANDROID]"
$ws.Range("C27").Value = "C70804"

# --- Move the saved selection to match the post-edit workbook state ---
$ws.Range("E28").Select() | Out-Null

Write-Output "done"
